# "Generate Report for handback"
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handoffs have now been handed back (in sync with en-US):
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: status text updated, and the "Latest Target File"
#    (E) / "Latest Handback File" (F) / "Latest Handback DateTime" (G)
#    columns are filled in for the two real rows (a.md.md, b.md.md)

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/797a4275f7a0dbb420b3b6fcbdd0ed1d1835dfcc/e2e/a.md.md"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18baad778107421f565db2094e5a4626adc2f07c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
$zhXlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $zhTargetUrl, "", "", "a.md.md")
$wsZh.Range("E2").Font.Underline = 2
$wsZh.Range("E2").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhHandbackUrl, "", "", $zhXlfName)
$wsZh.Range("F2").Font.Underline = 2
$wsZh.Range("F2").Font.Color = 15570276

$wsZh.Range("G2").Value = "2016-01-26 09:33:13"

$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $zhTargetUrl, "", "", "a.md.md")
$wsZh.Range("E3").Font.Underline = 2
$wsZh.Range("E3").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhHandbackUrl, "", "", $zhXlfName)
$wsZh.Range("F3").Font.Underline = 2
$wsZh.Range("F3").Font.Color = 15570276

$wsZh.Range("G3").Value = "2016-01-26 09:33:13"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/797a4275f7a0dbb420b3b6fcbdd0ed1d1835dfcc/e2e/a.md.md"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7d5f316ee658a412f28ef78d24dce7ba2c508244/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
$deXlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $deTargetUrl, "", "", "a.md.md")
$wsDe.Range("E2").Font.Underline = 2
$wsDe.Range("E2").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deHandbackUrl, "", "", $deXlfName)
$wsDe.Range("F2").Font.Underline = 2
$wsDe.Range("F2").Font.Color = 15570276

$wsDe.Range("G2").Value = "2016-01-26 09:33:32"

$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $deTargetUrl, "", "", "a.md.md")
$wsDe.Range("E3").Font.Underline = 2
$wsDe.Range("E3").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deHandbackUrl, "", "", $deXlfName)
$wsDe.Range("F3").Font.Underline = 2
$wsDe.Range("F3").Font.Color = 15570276

$wsDe.Range("G3").Value = "2016-01-26 09:33:32"
